$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.778.56'
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").Value = '3.787.02'
$ws.Range("E3").Value = '  -1.18%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = "'435.93"
$ws.Range("E5").Value = '  +1.57%  '

$ws.Range("D6").Value = "'139.23"
$ws.Range("E6").Value = '  +6.41%  '

$ws.Range("D7").Value = "'0.621"
$ws.Range("E7").Value = '  +1.64%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  +0.78%  '

$ws.Range("D10").Value = "'0.153"
$ws.Range("E10").Value = '  -8.68%  '

$ws.Range("D11").Value = "'0.0000314"
$ws.Range("E11").Value = '  -14.14%  '

$ws.Range("D12").Value = "'42.88"
$ws.Range("E12").Value = '  +5.09%  '

$ws.Range("D13").Value = "'10.38"
$ws.Range("E13").Value = '  +3.33%  '

$ws.Range("D14").Value = '4.379.26'
$ws.Range("E14").Value = '  -1.37%  '

$ws.Range("D15").Value = "'14.84"
$ws.Range("E15").Value = '  -5.08%  '

$ws.Range("E16").Value = '  -0.48%  '

$ws.Range("D17").Value = '3.763.91'
$ws.Range("E17").Value = '  -1.48%  '

$ws.Range("E18").Value = '  +1.61%  '

$ws.Range("E19").Value = '  +7.00%  '

$ws.Range("D20").Value = '66.761.86'
$ws.Range("E20").Value = '  -0.45%  '

$ws.Range("D21").Value = "'422.74"
$ws.Range("E21").Value = '  +3.31%  '

$ws.Range("E22").Value = '  +1.49%  '

$ws.Range("D23").Value = "'3.23"
$ws.Range("E23").Value = '  +6.75%  '

$ws.Range("D24").Value = "'85.90"
$ws.Range("E24").Value = '  +0.71%  '

$ws.Range("D25").Value = "'37.31"
$ws.Range("E25").Value = '  +1.24%  '

$ws.Range("D26").Value = "'3.38"
$ws.Range("E26").Value = '  +3.71%  '

$ws.Range("D27").Value = "'9.83"
$ws.Range("E27").Value = '  +36.59%  '

$ws.Range("D28").Value = "'5.56"
$ws.Range("E28").Value = '  -1.85%  '

$ws.Range("D29").Value = "'9.79"
$ws.Range("E29").Value = '  +3.23%  '

$ws.Range("D30").Value = "'734.15"
$ws.Range("E30").Value = '  +6.80%  '

$ws.Range("D31").Value = "'13.79"
$ws.Range("E31").Value = '  +10.97%  '

$ws.Range("E32").Value = '  +10.16%  '

$ws.Range("E33").Value = '  +2.64%  '

$ws.Range("D34").Value = "'42.88"
$ws.Range("E34").Value = '  +11.14%  '

$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("E36").Value = '  +1.84%  '

$ws.Range("D37").Value = "'5.60"
$ws.Range("E37").Value = '  +23.74%  '

$ws.Range("D38").Value = "'56.28"
$ws.Range("E38").Value = '  +2.09%  '

$ws.Range("D39").Value = "'0.0480"
$ws.Range("E39").Value = '  +5.03%  '

$ws.Range("E40").Value = '  +40.70%  '

$ws.Range("E41").Value = '  -3.61%  '

$ws.Range("E42").Value = '  +3.63%  '

$ws.Range("E43").Value = '  -0.07%  '

$ws.Range("D44").Value = '0.0₃0674'
$ws.Range("E44").Value = '  -14.89%  '

$ws.Range("D45").Value = "'0.335"
$ws.Range("E45").Value = '  +13.84%  '

$ws.Range("D46").Value = "'3.29"
$ws.Range("E46").Value = '  +5.58%  '

$ws.Range("B47").Value = 'LidoDAOToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D47").Value = "'3.31"
$ws.Range("E47").Value = '  +1.01%  '

$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = "'2.67"
$ws.Range("E48").Value = '  +5.25%  '

$ws.Range("E49").Value = '  +0.28%  '

$ws.Range("D50").Value = "'141.99"
$ws.Range("E50").Value = '  -4.28%  '

$ws.Range("D51").Value = "'2.82"
$ws.Range("E51").Value = '  +1.26%  '

# Reset style to remove quote-prefix formatting artifacts introduced above
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
